$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the AF2 related_works list with the new set of links
$ws.Range("AF2").Value = '''c("https://openalex.org/W4250942514", "https://openalex.org/W3160375675", "https://openalex.org/W4253897078", "https://openalex.org/W4283163498", "https://openalex.org/W3164536060", "https://openalex.org/W4221010187", "https://openalex.org/W4283257700", "https://openalex.org/W4226088185", "https://openalex.org/W4320037181", "https://openalex.org/W3005980542")'

# Insert two new blank rows at row 5, pushing the existing row 5
# record (and the rows below it) down to make room for two new
# article records sourced above it.
$ws.Rows.Item(5).Insert()
$ws.Rows.Item(6).Insert()

# Populate new row 5
$ws.Range("A5").Value = '''https://openalex.org/W4388703944'
$ws.Range("B5").Value = '''Bloch-Sulzberger Syndrome: A Rare X-Linked Dominant Genetic Disorder in a Newborn'
$ws.Range("C5").Value = '''list(au_id = c("https://openalex.org/A5093268654", "https://openalex.org/A5002855036", "https://openalex.org/A5041852419", "https://openalex.org/A5059557848", "https://openalex.org/A5022844184", "https://openalex.org/A5021789119", "https://openalex.org/A5092343508", "https://openalex.org/A5091408968"), au_display_name = c("Utsav P Vaghani", "Abdul Khalid Qadree", "S. C. Mehta", "Nalin Chaudhary", "Kusum Sharma", "Sanjeev Chaudhary", "Anasonye Emmanuel Kelechi", "Kausar Bano"), au_orcid = c(NA, NA, 
NA, NA, "https://orcid.org/0009-0001-7021-654X", NA, "https://orcid.org/0009-0002-9278-4713", NA), author_position = c("first", "middle", "middle", "middle", "middle", "middle", "middle", "last"), au_affiliation_raw = c("", "", "", "", "", "", "", ""), institution_id = c(NA, NA, NA, NA, NA, NA, NA, NA), institution_display_name = c(NA, NA, NA, NA, NA, NA, NA, NA), institution_ror = c(NA, NA, NA, NA, NA, NA, NA, NA), institution_country_code = c(NA, NA, NA, NA, NA, NA, NA, NA), institution_type = c(NA, 
NA, NA, NA, NA, NA, NA, NA), institution_lineage = c(NA, NA, NA, NA, NA, NA, NA, NA))'
$ws.Range("D5").Value = '''Bloch-Sulzberger Syndrome, also known as Incontinence Pigmentosa (IP), is a rare genodermatosis in which skin involvement occurs in almost all patients. Additionally, other ectodermal tissues like the central nervous system, eyes, hair, nails, and teeth may also be impacted. An X-linked dominant inheritance pattern characterizes the condition. But in our situation, IP caused a mutation in the body cells. There are four steps to the dermatological results. We describe the case of a 12-day-old female who had cutaneous features. It is crucial to make an early diagnosis using criteria like cutaneous symptoms so that quick diagnoses and interventions for other organs can be made to control more deadly complications in the future.'
$ws.Range("E5").Value = '''2023-11-15'
$ws.Range("F5").Value = '''Cureus'
$ws.Range("G5").Value = '''https://openalex.org/S2738950867'
$ws.Range("H5").Value = '''Cureus, Inc.'
$ws.Range("I5").Value = '''2168-8184'
$ws.Range("J5").Value = '''https://doi.org/10.7759/cureus.48823'
$ws.Range("K5").Value = '''https://assets.cureus.com/uploads/case_report/pdf/193461/20231115-17608-1mlie0r.pdf'
$ws.Range("L5").Value = '''N/A'
$ws.Range("M5").Value = '''publishedVersion'
$ws.Range("N5").Value = '''N/A'
$ws.Range("O5").Value = '''N/A'
$ws.Range("P5").Value = '''N/A'
$ws.Range("Q5").Value = '''N/A'
$ws.Range("R5").Value = '''TRUE'
$ws.Range("S5").Value = '''TRUE'
$ws.Range("T5").Value = '''bronze'
$ws.Range("U5").Value = '''https://assets.cureus.com/uploads/case_report/pdf/193461/20231115-17608-1mlie0r.pdf'
$ws.Range("V5").Value = '''FALSE'
$ws.Range("W5").Value = '''en'
$ws.Range("X5").Value = '''N/A'
$ws.Range("Y5").Value = '''0'
$ws.Range("Z5").Value = '''2023'
$ws.Range("AA5").Value = '''https://api.openalex.org/works?filter=cites:W4388703944'
$ws.Range("AB5").Value = '''c(openalex = "https://openalex.org/W4388703944", doi = "https://doi.org/10.7759/cureus.48823", pmid = "https://pubmed.ncbi.nlm.nih.gov/38106755")'
$ws.Range("AC5").Value = '''https://doi.org/10.7759/cureus.48823'
$ws.Range("AD5").Value = '''article'
$ws.Range("AE5").Value = '''c("https://openalex.org/W2053585674", "https://openalex.org/W2410764750", "https://openalex.org/W2753899911", "https://openalex.org/W4226417483", "https://openalex.org/W4241981241", "https://openalex.org/W4250539920")'
$ws.Range("AF5").Value = '''c("https://openalex.org/W4242331595", "https://openalex.org/W2114417421", "https://openalex.org/W2008290453", "https://openalex.org/W2412968475", "https://openalex.org/W2048373844", "https://openalex.org/W2057126451", "https://openalex.org/W2113956110", "https://openalex.org/W2029033789", "https://openalex.org/W1521115193", "https://openalex.org/W153882164")'
$ws.Range("AG5").Value = '''FALSE'
$ws.Range("AH5").Value = '''FALSE'

# Populate new row 6
$ws.Range("A6").Value = '''https://openalex.org/W4389160784'
$ws.Range("B6").Value = '''Futuristic Novel Therapeutic Approaches in the Treatment of Rheumatoid Arthritis'
$ws.Range("C6").Value = '''list(au_id = c("https://openalex.org/A5049453970", "https://openalex.org/A5014813333", "https://openalex.org/A5022343185", "https://openalex.org/A5021789119", "https://openalex.org/A5070161110", "https://openalex.org/A5093382818"), au_display_name = c("Vikramaditya Rai", "Nirmal Patel", "Simi R Mammen", "Sanjeev Chaudhary", "Sanan Arshad", "Shahzad W Munazzam"), au_orcid = c("https://orcid.org/0000-0002-2047-0531", "https://orcid.org/0000-0002-9799-7868", NA, NA, NA, NA), author_position = c("first", 
"middle", "middle", "middle", "middle", "last"), au_affiliation_raw = c("", "", "", "", "", ""), institution_id = c(NA, NA, NA, NA, NA, NA), institution_display_name = c(NA, NA, NA, NA, NA, NA), institution_ror = c(NA, NA, NA, NA, NA, NA), institution_country_code = c(NA, NA, NA, NA, NA, NA), institution_type = c(NA, NA, NA, NA, NA, NA), institution_lineage = c(NA, NA, NA, NA, NA, NA))'
$ws.Range("D6").Value = '''Rheumatoid arthritis (RA) is a chronic autoimmune disease characterized by systemic inflammation and joint destruction, leading to significant morbidity and reduced quality of life. Although significant progress has been made in the management of RA over the past few decades, many patients still fail to respond adequately to currently available therapies. This article aims to review the current landscape of RA treatment and explore potential novel therapeutic approaches that hold promise for the future. Advances in our understanding of the underlying pathogenesis of the disease have led to the identification of new targets and the development of innovative treatment strategies. This review focuses on emerging therapies including small molecule inhibitors, targeted biologics, cell-based therapies, and gene editing technologies that have shown potential in preclinical and early clinical trials. Additionally, we discuss the challenges and opportunities associated with the use of these new approaches in the treatment of RA. By elucidating the future of novel therapeutic approaches, this article provides insights that can guide clinicians and researchers in their efforts to improve outcomes for patients with RA.'
$ws.Range("E6").Value = '''2023-11-30'
$ws.Range("F6").Value = '''Cureus'
$ws.Range("G6").Value = '''https://openalex.org/S2738950867'
$ws.Range("H6").Value = '''Cureus, Inc.'
$ws.Range("I6").Value = '''2168-8184'
$ws.Range("J6").Value = '''https://doi.org/10.7759/cureus.49738'
$ws.Range("K6").Value = '''https://assets.cureus.com/uploads/review_article/pdf/197973/20231130-8005-1143v1s.pdf'
$ws.Range("L6").Value = '''N/A'
$ws.Range("M6").Value = '''publishedVersion'
$ws.Range("N6").Value = '''N/A'
$ws.Range("O6").Value = '''N/A'
$ws.Range("P6").Value = '''N/A'
$ws.Range("Q6").Value = '''N/A'
$ws.Range("R6").Value = '''TRUE'
$ws.Range("S6").Value = '''TRUE'
$ws.Range("T6").Value = '''bronze'
$ws.Range("U6").Value = '''https://assets.cureus.com/uploads/review_article/pdf/197973/20231130-8005-1143v1s.pdf'
$ws.Range("V6").Value = '''FALSE'
$ws.Range("W6").Value = '''en'
$ws.Range("X6").Value = '''N/A'
$ws.Range("Y6").Value = '''0'
$ws.Range("Z6").Value = '''2023'
$ws.Range("AA6").Value = '''https://api.openalex.org/works?filter=cites:W4389160784'
$ws.Range("AB6").Value = '''c(openalex = "https://openalex.org/W4389160784", doi = "https://doi.org/10.7759/cureus.49738", pmid = "https://pubmed.ncbi.nlm.nih.gov/38161868")'
$ws.Range("AC6").Value = '''https://doi.org/10.7759/cureus.49738'
$ws.Range("AD6").Value = '''article'
$ws.Range("AE6").Value = '''c("https://openalex.org/W2038593306", "https://openalex.org/W2046957298", "https://openalex.org/W2201326611", "https://openalex.org/W2597733778", "https://openalex.org/W2767394258", "https://openalex.org/W2803973097", "https://openalex.org/W2954501405", "https://openalex.org/W2976379921", "https://openalex.org/W2995734501", "https://openalex.org/W3180882235", "https://openalex.org/W3184467671", "https://openalex.org/W4200454701", "https://openalex.org/W4210366551", "https://openalex.org/W4226275346", 
"https://openalex.org/W4280590568", "https://openalex.org/W4282934357", "https://openalex.org/W4307106007", "https://openalex.org/W4311307625", "https://openalex.org/W4311479028", "https://openalex.org/W4315433173", "https://openalex.org/W4321166025", "https://openalex.org/W4323652643", "https://openalex.org/W4377164070", "https://openalex.org/W4383710632")'
$ws.Range("AF6").Value = '''c("https://openalex.org/W1967563157", "https://openalex.org/W1501007924", "https://openalex.org/W1143447605", "https://openalex.org/W2111650087", "https://openalex.org/W1995085194", "https://openalex.org/W2185550546", "https://openalex.org/W2416097977", "https://openalex.org/W2993678931", "https://openalex.org/W2028766184", "https://openalex.org/W3045583548")'
$ws.Range("AG6").Value = '''FALSE'
$ws.Range("AH6").Value = '''FALSE'

